$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 20833912
$ws.Range("I33").Value = 581.3333
$ws.Range("J33").Value = 41667244
$ws.Range("K33").Value = 581.3333
$ws.Range("L33").Value = 41667244
$ws.Range("M33").Value = -352.3333
$ws.Range("N33").Value = -41667702

$ws.Range("H40").Value = 2943.3
$ws.Range("I40").Value = 2821.6667
$ws.Range("J40").Value = 2995.4285
$ws.Range("K40").Value = 2821.6667
$ws.Range("L40").Value = 2995.4285
$ws.Range("M40").Value = -2646.6667
$ws.Range("N40").Value = -3345.4285

$ws.Range("H80").Value = 947696.4399999999
$ws.Range("I80").Value = 1748872
$ws.Range("K80").Value = 5246616
$ws.Range("M80").Value = -5245618

$ws.Range("H83").Value = 947696.4399999999
$ws.Range("I83").Value = 1748872
$ws.Range("K83").Value = 15739848
$ws.Range("M83").Value = -15734856

$ws.Range("H104").Value = 1868.6666
$ws.Range("I104").Value = 1062.25
$ws.Range("J104").Value = 3481.5
$ws.Range("K104").Value = 3186.75
$ws.Range("L104").Value = 10444.5
$ws.Range("M104").Value = -1439.75
$ws.Range("N104").Value = -13938.5

$ws.Range("H106").Value = 5470.3335
$ws.Range("I106").Value = 2020.75
$ws.Range("K106").Value = 2020.75
$ws.Range("M106").Value = -1389.75

$ws.Range("H123").Value = 86159.71000000001
$ws.Range("J123").Value = 88853
$ws.Range("L123").Value = 88853
$ws.Range("N123").Value = -98653

$ws.Range("H132").Value = 3061.4666
$ws.Range("I132").Value = 2918.0754
$ws.Range("K132").Value = 8754.226200000001
$ws.Range("M132").Value = -6224.226200000001

$ws.Range("H135").Value = 1041.1364
$ws.Range("I135").Value = 1041.1364
$ws.Range("K135").Value = 9370.2276
$ws.Range("M135").Value = -6835.2276

$ws.Range("H137").Value = 2337.11
$ws.Range("I137").Value = 1631.7693
$ws.Range("J137").Value = 2442.5059
$ws.Range("K137").Value = 4895.3079
$ws.Range("L137").Value = 7327.5177
$ws.Range("M137").Value = -2345.3079
$ws.Range("N137").Value = -12427.5177

$ws.Range("H138").Value = 3116.772
$ws.Range("I138").Value = 2186.08
$ws.Range("J138").Value = 3843.875
$ws.Range("K138").Value = 6558.24
$ws.Range("L138").Value = 11531.625
$ws.Range("M138").Value = -1418.24
$ws.Range("N138").Value = -21811.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 6725
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 8060
$ws.Range("K5").Value = 50
$ws.Range("L5").Value = 8060
$ws.Range("M5").Value = 62
$ws.Range("N5").Value = -8284

$ws.Range("H32").Value = 547.0599999999999
$ws.Range("I32").Value = 335.60214
$ws.Range("K32").Value = 335.60214
$ws.Range("M32").Value = -48.60214000000002

$ws.Range("H45").Value = 3273
$ws.Range("I45").Value = 3129
$ws.Range("K45").Value = 3129
$ws.Range("M45").Value = -2752

$ws.Range("H61").Value = 2501495
$ws.Range("I61").Value = 1668160
$ws.Range("K61").Value = 1668160
$ws.Range("M61").Value = -1667948

$ws.Range("H97").Value = 274.1875
$ws.Range("I97").Value = 229.78572
$ws.Range("K97").Value = 229.78572
$ws.Range("M97").Value = 266.21428

$ws.Range("H122").Value = 3946.8462
$ws.Range("J122").Value = 4416.6665
$ws.Range("L122").Value = 13249.9995
$ws.Range("N122").Value = -18149.9995

$ws.Range("H132").Value = 1178517.1
$ws.Range("I132").Value = 834857.7
$ws.Range("J132").Value = 2003300
$ws.Range("K132").Value = 2504573.1
$ws.Range("L132").Value = 6009900
$ws.Range("M132").Value = -2502043.1
$ws.Range("N132").Value = -6014960

$ws.Range("H136").Value = 2501495
$ws.Range("I136").Value = 1668160
$ws.Range("K136").Value = 5004480
$ws.Range("M136").Value = -5001930

$ws.Range("H141").Value = 59999
$ws.Range("J141").Value = 59999
$ws.Range("L141").Value = 59999
$ws.Range("N141").Value = -70359

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 6725
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 8060
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 8060
$ws.Range("M4").Value = 65
$ws.Range("N4").Value = -8290

$ws.Range("H20").Value = 1861.3143
$ws.Range("I20").Value = 1808.5834
$ws.Range("J20").Value = 1976.3636
$ws.Range("K20").Value = 1808.5834
$ws.Range("L20").Value = 1976.3636
$ws.Range("M20").Value = -1561.5834
$ws.Range("N20").Value = -2470.3636

$ws.Range("H22").Value = 416.83334
$ws.Range("I22").Value = 420.2
$ws.Range("K22").Value = 420.2
$ws.Range("M22").Value = -247.2

$ws.Range("H86").Value = 1706.3572
$ws.Range("I86").Value = 1762.6364
$ws.Range("K86").Value = 1762.6364
$ws.Range("M86").Value = -639.6364000000001

$ws.Range("H89").Value = 1706.3572
$ws.Range("I89").Value = 1762.6364
$ws.Range("K89").Value = 8813.182000000001
$ws.Range("M89").Value = -3197.182000000001

$ws.Range("H107").Value = 13526165
$ws.Range("J107").Value = 41686116
$ws.Range("L107").Value = 41686116
$ws.Range("N107").Value = -41689956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2481.4
$ws.Range("I31").Value = 729.7895
$ws.Range("J31").Value = 2892.2715
$ws.Range("K31").Value = 729.7895
$ws.Range("L31").Value = 2892.2715
$ws.Range("M31").Value = -434.7895
$ws.Range("N31").Value = -3482.2715

$ws.Range("H34").Value = 2481.4
$ws.Range("I34").Value = 729.7895
$ws.Range("J34").Value = 2892.2715
$ws.Range("K34").Value = 729.7895
$ws.Range("L34").Value = 2892.2715
$ws.Range("M34").Value = -527.7895
$ws.Range("N34").Value = -3296.2715

$ws.Range("H132").Value = 2852
$ws.Range("I132").Value = 2548.75
$ws.Range("K132").Value = 7646.25
$ws.Range("M132").Value = -5116.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1298.8
$ws.Range("I132").Value = 998
$ws.Range("J132").Value = 1374
$ws.Range("K132").Value = 8982
$ws.Range("L132").Value = 12366
$ws.Range("M132").Value = -6452
$ws.Range("N132").Value = -17426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 3000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 3000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 3000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -3780

$ws.Range("H80").Value = 3600
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3600
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3600
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -5596

$ws.Range("H83").Value = 3600
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3600
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 18000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -27984

$ws.Range("H97").Value = 1385.9286
$ws.Range("I97").Value = 1507.25
$ws.Range("J97").Value = 1082.625
$ws.Range("K97").Value = 1507.25
$ws.Range("L97").Value = 1082.625
$ws.Range("M97").Value = -1011.25
$ws.Range("N97").Value = -2074.625

$ws.Range("H113").Value = 4975
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H126").Value = 11894.571
$ws.Range("I126").Value = 17253
$ws.Range("J126").Value = 4750
$ws.Range("K126").Value = 51759
$ws.Range("L126").Value = 14250
$ws.Range("M126").Value = -49289
$ws.Range("N126").Value = -19190

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14292106
$ws.Range("I81").Value = 3966
$ws.Range("J81").Value = 18188872
$ws.Range("K81").Value = 7932
$ws.Range("L81").Value = 36377744
$ws.Range("M81").Value = -6871
$ws.Range("N81").Value = -36379866

$ws.Range("H84").Value = 14292106
$ws.Range("I84").Value = 3966
$ws.Range("J84").Value = 18188872
$ws.Range("K84").Value = 39660
$ws.Range("L84").Value = 181888720
$ws.Range("M84").Value = -34356
$ws.Range("N84").Value = -181899328

$ws.Range("H96").Value = 4369.357
$ws.Range("J96").Value = 4644.6665
$ws.Range("L96").Value = 4644.6665
$ws.Range("N96").Value = -7390.6665

$ws.Range("H126").Value = 1686.4286
$ws.Range("I126").Value = 1635.1333
$ws.Range("J126").Value = 1814.6666
$ws.Range("K126").Value = 4905.3999
$ws.Range("L126").Value = 5443.9998
$ws.Range("M126").Value = -2435.3999
$ws.Range("N126").Value = -10383.9998

$ws.Range("H132").Value = 694281.4399999999
$ws.Range("I132").Value = 1181246.9
$ws.Range("J132").Value = 4413.6665
$ws.Range("K132").Value = 3543740.7
$ws.Range("L132").Value = 13240.9995
$ws.Range("M132").Value = -3541210.7
$ws.Range("N132").Value = -18300.9995
